$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 'AUB'
$ws.Cells.Item(3, 2).Value = 'AGR'
$ws.Cells.Item(4, 2).Value = 'ALG'
$ws.Cells.Item(5, 2).ClearFormats()
$ws.Cells.Item(5, 2).Value = 'ALN'
$ws.Cells.Item(6, 2).Value = 'ALH'
$ws.Cells.Item(7, 2).Value = 'ACT'
$ws.Cells.Item(8, 2).ClearFormats()
$ws.Cells.Item(8, 2).Value = 'AJA'
$ws.Cells.Item(9, 2).Value = 'ALD'
$ws.Cells.Item(10, 2).Value = 'APR'
$ws.Cells.Item(11, 2).Value = 'APAC'
$ws.Cells.Item(12, 2).Value = 'ARG'
$ws.Cells.Item(13, 2).Value = 'ARR'
$ws.Cells.Item(14, 2).Value = 'ARU'
$ws.Cells.Item(15, 2).Value = 'ARE'
$ws.Cells.Item(16, 2).Value = 'ABN'
$ws.Cells.Item(17, 2).Value = 'AREL'
$ws.Cells.Item(18, 2).ClearFormats()
$ws.Cells.Item(18, 2).Value = 'ARS'
$ws.Cells.Item(19, 2).Value = 'ASS'
$ws.Cells.Item(20, 2).Value = 'BTÇ'
$ws.Cells.Item(21, 2).Value = 'BNN'
$ws.Cells.Item(22, 2).Value = 'BRU'
$ws.Cells.Item(23, 2).Value = 'BSR'
$ws.Cells.Item(24, 2).Value = 'BDS'
$ws.Cells.Item(25, 2).ClearFormats()
$ws.Cells.Item(25, 2).Value = 'BSM'
$ws.Cells.Item(26, 2).Value = 'BYX'
$ws.Cells.Item(27, 2).Value = 'BLM'
$ws.Cells.Item(28, 2).Value = 'BBC'
$ws.Cells.Item(29, 2).Value = 'BBT'
$ws.Cells.Item(30, 2).Value = 'BVA'
$ws.Cells.Item(31, 2).Value = 'BVA'
$ws.Cells.Item(32, 2).Value = 'BMJ'
$ws.Cells.Item(33, 2).Value = 'BSC'
$ws.Cells.Item(34, 2).Value = 'BSF'
$ws.Cells.Item(35, 2).Value = 'BQO'
$ws.Cells.Item(36, 2).Value = 'BBM'
$ws.Cells.Item(37, 2).Value = 'BDC'
$ws.Cells.Item(38, 2).ClearFormats()
$ws.Cells.Item(38, 2).Value = 'BDS'
$ws.Cells.Item(39, 2).Value = 'CAP'
$ws.Cells.Item(40, 2).Value = 'CBC'
$ws.Cells.Item(41, 2).Value = 'CBD'
$ws.Cells.Item(42, 2).Value = 'CDI'
$ws.Cells.Item(43, 2).ClearFormats()
$ws.Cells.Item(43, 2).Value = 'CDA'
$ws.Cells.Item(44, 2).ClearFormats()
$ws.Cells.Item(44, 2).Value = 'CDD'
$ws.Cells.Item(45, 2).Value = 'CCB'
$ws.Cells.Item(46, 2).Value = 'CCR'
$ws.Cells.Item(47, 2).Value = 'CJZ'
$ws.Cells.Item(48, 2).Value = 'CZR'
$ws.Cells.Item(49, 2).Value = 'CBA'
$ws.Cells.Item(50, 2).Value = 'CML'
$ws.Cells.Item(51, 2).Value = 'CGA'
$ws.Cells.Item(52, 2).Value = 'CPM'
$ws.Cells.Item(53, 2).Value = 'CUB'
$ws.Cells.Item(54, 2).Value = 'CRR'
$ws.Cells.Item(55, 2).Value = 'CSG'
$ws.Cells.Item(56, 2).Value = 'CAT'
$ws.Cells.Item(57, 2).Value = 'CDR'
$ws.Cells.Item(58, 2).Value = 'CUT'
$ws.Cells.Item(59, 2).Value = 'CCÇ'
$ws.Cells.Item(60, 2).Value = 'CDD'
$ws.Cells.Item(61, 2).Value = 'CDN'
$ws.Cells.Item(62, 2).Value = 'CNG'
$ws.Cells.Item(63, 2).Value = 'CRM'
$ws.Cells.Item(64, 2).Value = 'CXX'
$ws.Cells.Item(65, 2).Value = 'CES'
$ws.Cells.Item(66, 2).Value = 'CBT'
$ws.Cells.Item(67, 2).Value = 'CTE'
$ws.Cells.Item(68, 2).Value = 'CDM'
$ws.Cells.Item(69, 2).Value = 'CIG'
$ws.Cells.Item(70, 2).Value = 'CDC'
$ws.Cells.Item(71, 2).Value = 'CLV'
$ws.Cells.Item(72, 2).Value = 'DMÃ'
$ws.Cells.Item(73, 2).Value = 'DET'
$ws.Cells.Item(74, 2).Value = 'DIT'
$ws.Cells.Item(75, 2).Value = 'DIN'
$ws.Cells.Item(76, 2).Value = 'DUE'
$ws.Cells.Item(77, 2).Value = 'EMA'
$ws.Cells.Item(78, 2).ClearFormats()
$ws.Cells.Item(78, 2).Value = 'ESP'
$ws.Cells.Item(79, 2).Value = 'FUD'
$ws.Cells.Item(80, 2).Value = 'FMO'
$ws.Cells.Item(81, 2).ClearFormats()
$ws.Cells.Item(81, 2).Value = 'GBV'
$ws.Cells.Item(82, 2).Value = 'GRB'
$ws.Cells.Item(83, 2).Value = 'GRM'
$ws.Cells.Item(84, 2).Value = 'GJÃ'
$ws.Cells.Item(85, 2).Value = 'IBI'
$ws.Cells.Item(86, 2).Value = 'IGY'
$ws.Cells.Item(87, 2).Value = 'IMC'
$ws.Cells.Item(88, 2).Value = 'ING'
$ws.Cells.Item(89, 2).Value = 'ITB'
$ws.Cells.Item(90, 2).Value = 'ITG'
$ws.Cells.Item(91, 2).Value = 'ITR'
$ws.Cells.Item(92, 2).Value = 'TIB'
$ws.Cells.Item(93, 2).Value = 'JCR'
$ws.Cells.Item(94, 2).Value = 'JRC'
$ws.Cells.Item(95, 2).Value = 'JPA'
$ws.Cells.Item(96, 2).Value = 'JCD'
$ws.Cells.Item(97, 2).Value = 'JTR'
$ws.Cells.Item(98, 2).Value = 'JZN'
$ws.Cells.Item(99, 2).Value = 'JDS'
$ws.Cells.Item(100, 2).Value = 'JPG'
$ws.Cells.Item(101, 2).Value = 'JUR'
$ws.Cells.Item(102, 2).Value = 'LGO'
$ws.Cells.Item(103, 2).Value = 'LDD'
$ws.Cells.Item(104, 2).Value = 'GS'
$ws.Cells.Item(105, 2).Value = 'LAT'
$ws.Cells.Item(106, 2).Value = 'LIVR'
$ws.Cells.Item(107, 2).Value = 'LGR'
$ws.Cells.Item(108, 2).Value = 'LUC'
$ws.Cells.Item(109, 2).Value = 'MDA'
$ws.Cells.Item(110, 2).Value = 'MAT'
$ws.Cells.Item(111, 2).Value = 'MGP'
$ws.Cells.Item(112, 2).ClearFormats()
$ws.Cells.Item(112, 2).Value = 'MAÍ'
$ws.Cells.Item(113, 2).Value = 'MAC'
$ws.Cells.Item(114, 2).Value = 'MAR'
$ws.Cells.Item(115, 2).Value = 'MZP'
$ws.Cells.Item(116, 2).Value = 'MSS'
$ws.Cells.Item(117, 2).Value = 'MTR'
$ws.Cells.Item(118, 2).Value = 'MTH'
$ws.Cells.Item(119, 2).Value = 'MTG'
$ws.Cells.Item(120, 2).Value = 'MTE'
$ws.Cells.Item(121, 2).ClearFormats()
$ws.Cells.Item(121, 2).Value = 'MGR'
$ws.Cells.Item(122, 2).Value = 'MTD'
$ws.Cells.Item(123, 2).ClearFormats()
$ws.Cells.Item(123, 2).Value = 'MTE'
$ws.Cells.Item(124, 2).Value = 'MON'
$ws.Cells.Item(125, 2).Value = 'MUG'
$ws.Cells.Item(126, 2).Value = 'NAT'
$ws.Cells.Item(127, 2).Value = 'NZZ'
$ws.Cells.Item(128, 2).Value = 'NFT'
$ws.Cells.Item(129, 2).Value = 'NOL'
$ws.Cells.Item(130, 2).Value = 'NPL'
$ws.Cells.Item(131, 2).Value = 'ODA'
$ws.Cells.Item(132, 2).Value = 'OLD'
$ws.Cells.Item(133, 2).Value = 'OUV'
$ws.Cells.Item(134, 2).Value = 'PAT'
$ws.Cells.Item(135, 2).Value = 'PSS'
$ws.Cells.Item(136, 2).Value = 'PAT'
$ws.Cells.Item(137, 2).Value = 'PUL'
$ws.Cells.Item(138, 2).Value = 'PBA'
$ws.Cells.Item(139, 2).Value = 'PLA'
$ws.Cells.Item(140, 2).Value = 'PDF'
$ws.Cells.Item(141, 2).Value = 'POR'
$ws.Cells.Item(142, 2).Value = 'PNC'
$ws.Cells.Item(143, 2).Value = 'PCI'
$ws.Cells.Item(144, 2).Value = 'PIL'
$ws.Cells.Item(145, 2).Value = 'POS'
$ws.Cells.Item(146, 2).Value = 'PLZ'
$ws.Cells.Item(147, 2).Value = 'PPT'
$ws.Cells.Item(148, 2).Value = 'PIT'
$ws.Cells.Item(149, 2).Value = 'PCH'
$ws.Cells.Item(150, 2).Value = 'PDT'
$ws.Cells.Item(151, 2).Value = 'PJM'
$ws.Cells.Item(152, 2).Value = 'PBL'
$ws.Cells.Item(153, 2).Value = 'PRT'
$ws.Cells.Item(154, 2).Value = 'PRI'
$ws.Cells.Item(155, 2).Value = 'PNN'
$ws.Cells.Item(156, 2).Value = 'QMD'
$ws.Cells.Item(157, 2).Value = 'QXB'
$ws.Cells.Item(158, 2).Value = 'RMG'
$ws.Cells.Item(159, 2).Value = 'RAÇ'
$ws.Cells.Item(160, 2).Value = 'RDB'
$ws.Cells.Item(161, 2).Value = 'RDP'
$ws.Cells.Item(162, 2).ClearFormats()
$ws.Cells.Item(162, 2).Value = 'RSA'
$ws.Cells.Item(163, 2).Value = 'RDC'
$ws.Cells.Item(164, 2).Value = 'RTO'
$ws.Cells.Item(165, 2).Value = 'SGD'
$ws.Cells.Item(166, 2).Value = 'SDF'
$ws.Cells.Item(167, 2).Value = 'ACE'
$ws.Cells.Item(168, 2).Value = 'STC'
$ws.Cells.Item(169, 2).Value = 'STH'
$ws.Cells.Item(170, 2).Value = 'STI'
$ws.Cells.Item(171, 2).Value = 'ATL'
$ws.Cells.Item(172, 2).Value = 'STR'
$ws.Cells.Item(173, 2).Value = 'STT'
$ws.Cells.Item(174, 2).Value = 'STM'
$ws.Cells.Item(175, 2).Value = 'STG'
$ws.Cells.Item(176, 2).Value = 'ATA'
$ws.Cells.Item(177, 2).Value = 'SAB'
$ws.Cells.Item(178, 2).Value = 'SBT'
$ws.Cells.Item(179, 2).Value = 'SDG'
$ws.Cells.Item(180, 2).Value = 'SDC'
$ws.Cells.Item(181, 2).Value = 'SFC'
$ws.Cells.Item(182, 2).Value = 'SJC'
$ws.Cells.Item(183, 2).Value = 'SJR'
$ws.Cells.Item(184, 2).Value = 'SJT'
$ws.Cells.Item(185, 2).Value = 'SJL'
$ws.Cells.Item(186, 2).Value = 'SJC'
$ws.Cells.Item(187, 2).Value = 'SJE'
$ws.Cells.Item(188, 2).Value = 'SJP'
$ws.Cells.Item(189, 2).Value = 'SJP'
$ws.Cells.Item(190, 2).Value = 'SJB'
$ws.Cells.Item(191, 2).Value = 'SJB'
$ws.Cells.Item(192, 2).Value = 'SJS'
$ws.Cells.Item(193, 2).Value = 'SJC'
$ws.Cells.Item(194, 2).Value = 'SJR'
$ws.Cells.Item(195, 2).Value = 'SMM'
$ws.Cells.Item(196, 2).Value = 'SMT'
$ws.Cells.Item(197, 2).Value = 'SSR'
$ws.Cells.Item(198, 2).Value = 'SSU'
$ws.Cells.Item(199, 2).Value = 'SVS'
$ws.Cells.Item(200, 2).Value = 'SAP'
$ws.Cells.Item(201, 2).Value = 'SEB'
$ws.Cells.Item(202, 2).Value = 'SDR'
$ws.Cells.Item(203, 2).Value = 'SRG'
$ws.Cells.Item(204, 2).Value = 'SER'
$ws.Cells.Item(205, 2).Value = 'SRR'
$ws.Cells.Item(206, 2).Value = 'STZ'
$ws.Cells.Item(207, 2).Value = 'SBD'
$ws.Cells.Item(208, 2).ClearFormats()
$ws.Cells.Item(208, 2).Value = 'SLN'
$ws.Cells.Item(209, 2).Value = 'SOL'
$ws.Cells.Item(210, 2).Value = 'SOS'
$ws.Cells.Item(211, 2).Value = 'SOU'
$ws.Cells.Item(212, 2).Value = 'SMÉ'
$ws.Cells.Item(213, 2).Value = 'TCM'
$ws.Cells.Item(214, 2).Value = 'TPA'
$ws.Cells.Item(215, 2).Value = 'TVR'
$ws.Cells.Item(216, 2).Value = 'TXE'
$ws.Cells.Item(217, 2).Value = 'TNR'
$ws.Cells.Item(218, 2).ClearFormats()
$ws.Cells.Item(218, 2).Value = 'TRI'
$ws.Cells.Item(219, 2).Value = 'URU'
$ws.Cells.Item(220, 2).Value = 'UMB'
$ws.Cells.Item(221, 2).Value = 'VZA'
$ws.Cells.Item(222, 2).Value = 'VPL'
$ws.Cells.Item(223, 2).Value = 'VSR'
$ws.Cells.Item(224, 2).Value = 'ZBL'
$ws.Cells.Item(1, 2).Value = 'sigla'

[void]$ws.Range("B4").Select()
